$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18-21: CV numbers column (D) changes from "Centers" to "-"
$ws.Range("D18").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("D20").Value = "-"
$ws.Range("D21").Value = "-"

# Row 22: new RF submission with all data
$ws.Range("A22").Value = "2023-03-04-2149_RF_all.csv"
$ws.Range("B22").Value = "RandomForest"
$ws.Range("C22").Value = "MoCo"
$ws.Range("D22").Value = "-"
$ws.Range("E22").Value = "-"
$ws.Range("F22").Value = "weakly supervision"
$ws.Range("M22").Value = "March 4, 2023, 8:50 p.m."
$ws.Range("N22").Value = 0.603

# Row 21: add Hand in / Test AUC values
$ws.Range("M21").Value = "March 4, 2023, 8:35 p.m."
$ws.Range("N21").Value = 0.647

# Expand the table ("Tabelle1") to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N22"))

# Scroll the view and select the newly added cell (mirrors the author's
# final on-screen state after adding the new submission row)
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M22").Select()
